$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the existing
# PercActivations/PercSegmentAreas/RelativeCAMImportance/PercActivationsRescaled
# columns (B:E) one place to the right (C:F), and the segment-name values
# that lived in column A stay in column A for now - we move them to the
# new column B below.
$ws.Columns.Item(2).Insert()

# New column header for the segment names that will move into column B.
$ws.Range("B1").Value = "segments"

# The new B1 cell was created blank by the column insert (there was nothing
# in the old B1 to inherit a style from), so it needs the same bold /
# centered / bordered header style the other header cells (C1:F1) carry.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$lastRow = 20

for ($r = 2; $r -le $lastRow; $r++) {
    # Segment name currently still in column A -> move to new column B.
    $name = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $name

    # Column A becomes a 0-based numeric segment index.
    $ws.Cells.Item($r, 1).Value = $r - 2

    # Segment-name cells in column B carry no explicit style (plain cells),
    # matching the rest of the data columns, while column A keeps the
    # original bordered/bold/centered header style.
    $ws.Cells.Item($r, 2).Style = "Normal"
}
